# In all models change "Cross references" to "Database references" to
# replace an ambiguous name with a clear one.

$wb = $excel.ActiveWorkbook

# Rename the "Cross references" sheet to "Database references".
$ws = $wb.Worksheets.Item("Cross references")
$ws.Name = "Database references"

# Make it the active/selected tab (was previously on "Biomass reactions").
$ws.Select()

# Turn on iterative calculation's max-change tolerance (1E-4).
$excel.MaxChange = 0.0001
